# Update the cached "today" text shown by the automatic date/time
# placeholders (type datetimeFigureOut) from 9/27/18 to 8/27/19 across
# every part of the deck that carries one: the slide master, every
# slide layout, and the notes master.

$p = $ppt.ActivePresentation

$oldDate = "9/27/18"
$newDate = "8/27/19"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }

        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }

        if ($isDatePlaceholder -and ($sh.TextFrame.TextRange.Text -eq $oldDate)) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout hanging off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholders $notesMaster.Shapes
